$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.014886991352967
$ws.Range("D2").Value = 1.041572210879171
$ws.Range("E2").Value = 1.01657958546885
$ws.Range("F2").Value = 1.044029674126431
$ws.Range("I2").Value = 1.035502744679201
$ws.Range("J2").Value = 1.020115658763645
$ws.Range("K2").Value = 1.044351346431286
$ws.Range("L2").Value = 1.019431334440237
$ws.Range("M2").Value = 1.046801873561701
$ws.Range("N2").Value = 1.010697558652584
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.016418288807846
$ws.Range("D3").Value = 1.042318752867021
$ws.Range("E3").Value = 1.017895117441466
$ws.Range("F3").Value = 1.045128995407642
$ws.Range("I3").Value = 1.035681005346287
$ws.Range("J3").Value = 1.021278658044432
$ws.Range("K3").Value = 1.044908711781824
$ws.Range("L3").Value = 1.020550865712495
$ws.Range("M3").Value = 1.047711597531805
$ws.Range("N3").Value = 1.011092836752121
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.017407435628509
$ws.Range("D4").Value = 1.042798975326388
$ws.Range("E4").Value = 1.018745178888278
$ws.Range("F4").Value = 1.045837340075268
$ws.Range("I4").Value = 1.035793280280198
$ws.Range("J4").Value = 1.022029201167756
$ws.Range("K4").Value = 1.045265866839377
$ws.Range("L4").Value = 1.021273576977102
$ws.Range("M4").Value = 1.048296679083774
$ws.Range("N4").Value = 1.011347632672596
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.017822873231924
$ws.Range("D5").Value = 1.043000181102251
$ws.Range("E5").Value = 1.019102269420221
$ws.Range("F5").Value = 1.046134414559089
$ws.Range("I5").Value = 1.035839745126762
$ws.Range("J5").Value = 1.022344257796123
$ws.Range("K5").Value = 1.045415177729779
$ws.Range("L5").Value = 1.021577003325444
$ws.Range("M5").Value = 1.048541794290335
$ws.Range("N5").Value = 1.011454517335175
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.017892603838573
$ws.Range("D6").Value = 1.043033924580831
$ws.Range("E6").Value = 1.019162210526257
$ws.Range("F6").Value = 1.046184252885758
$ws.Range("I6").Value = 1.035847503640817
$ws.Range("J6").Value = 1.022397129701411
$ws.Range("K6").Value = 1.045440198603146
$ws.Range("L6").Value = 1.021627926563867
$ws.Range("M6").Value = 1.048582900236595
$ws.Range("N6").Value = 1.011472450219768
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.017412988286335
$ws.Range("D7").Value = 1.042801666518722
$ws.Range("E7").Value = 1.018749951423627
$ws.Range("F7").Value = 1.04584131240067
$ws.Range("I7").Value = 1.035793904036246
$ws.Range("J7").Value = 1.022033412816041
$ws.Range("K7").Value = 1.045267865227926
$ws.Range("L7").Value = 1.021277632945469
$ws.Range("M7").Value = 1.04829995767524
$ws.Range("N7").Value = 1.011349061778069
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.01540485864976
$ws.Range("D8").Value = 1.041825096438466
$ws.Range("E8").Value = 1.017024422822764
$ws.Range("F8").Value = 1.044401815243998
$ws.Range("I8").Value = 1.035563624612976
$ws.Range("J8").Value = 1.020509116527412
$ws.Range("K8").Value = 1.044540435022315
$ws.Range("L8").Value = 1.019810040721867
$ws.Range("M8").Value = 1.047110059448721
$ws.Range("N8").Value = 1.01083134790601
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.011852837783804
$ws.Range("D9").Value = 1.040082490765393
$ws.Range("E9").Value = 1.01397453833955
$ws.Range("F9").Value = 1.041842240062579
$ws.Range("I9").Value = 1.035134323368208
$ws.Range("J9").Value = 1.01780756247623
$ws.Range("K9").Value = 1.043231808194725
$ws.Range("L9").Value = 1.017210680746407
$ws.Range("M9").Value = 1.044985876746627
$ws.Range("N9").Value = 1.009911509748014
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.009475251103775
$ws.Range("D10").Value = 1.038906094647093
$ws.Range("E10").Value = 1.011934648686264
$ws.Range("F10").Value = 1.040120246487484
$ws.Range("I10").Value = 1.03483231227235
$ws.Range("J10").Value = 1.015995692969898
$ws.Range("K10").Value = 1.042341342119609
$ws.Range("L10").Value = 1.015468492515486
$ws.Range("M10").Value = 1.043551188210768
$ws.Range("N10").Value = 1.009293074821675
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.00844333896347
$ws.Range("D11").Value = 1.038393216961082
$ws.Range("E11").Value = 1.011049692362291
$ws.Range("F11").Value = 1.039370866221991
$ws.Range("I11").Value = 1.034697787318662
$ws.Range("J11").Value = 1.015208478538605
$ws.Range("K11").Value = 1.041951471078047
$ws.Range("L11").Value = 1.014711824102829
$ws.Range("M11").Value = 1.042925518233904
$ws.Range("N11").Value = 1.009024021502527
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.00805966959656
$ws.Range("D12").Value = 1.038202186001902
$ws.Range("E12").Value = 1.010720721482819
$ws.Range("F12").Value = 1.039091946963307
$ws.Range("I12").Value = 1.034647254812343
$ws.Range("J12").Value = 1.014915664497793
$ws.Range("K12").Value = 1.041806009607585
$ws.Range("L12").Value = 1.014430412793468
$ws.Range("M12").Value = 1.042692446570679
$ws.Range("N12").Value = 1.008923890273573
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.008141984947186
$ws.Range("D13").Value = 1.038243186580426
$ws.Range("E13").Value = 1.010791298684306
$ws.Range("F13").Value = 1.039151801760725
$ws.Range("I13").Value = 1.034658119730445
$ws.Range("J13").Value = 1.014978492666759
$ws.Range("K13").Value = 1.041837240836341
$ws.Range("L13").Value = 1.014490792470078
$ws.Range("M13").Value = 1.042742471598224
$ws.Range("N13").Value = 1.008945377529334
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.008411632368247
$ws.Range("D14").Value = 1.038377437005192
$ws.Range("E14").Value = 1.01102250484897
$ws.Range("F14").Value = 1.039347822236047
$ws.Range("I14").Value = 1.034693621791782
$ws.Range("J14").Value = 1.015184282796032
$ws.Range("K14").Value = 1.041939460371769
$ws.Range("L14").Value = 1.014688569756718
$ws.Range("M14").Value = 1.042906266141673
$ws.Range("N14").Value = 1.009015748572516
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.00857772153453
$ws.Range("D15").Value = 1.038460083478798
$ws.Range("E15").Value = 1.011164924086014
$ws.Range("F15").Value = 1.039468521810998
$ws.Range("I15").Value = 1.034715421047249
$ws.Range("J15").Value = 1.015311022657482
$ws.Range("K15").Value = 1.042002355568624
$ws.Range("L15").Value = 1.014810380152521
$ws.Range("M15").Value = 1.043007096500466
$ws.Range("N15").Value = 1.009059080865253
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.009543684227805
$ws.Range("D16").Value = 1.038940059031106
$ws.Range("E16").Value = 1.01199334446871
$ws.Range("F16").Value = 1.040169901189431
$ws.Range("I16").Value = 1.034841161178141
$ws.Range("J16").Value = 1.016047881057537
$ws.Range("K16").Value = 1.042367126038583
$ws.Range("L16").Value = 1.015518661273858
$ws.Range("M16").Value = 1.043592618070785
$ws.Range("N16").Value = 1.009310904098542
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.010148956819027
$ws.Range("D17").Value = 1.039240199880032
$ws.Range("E17").Value = 1.012512537995553
$ws.Range("F17").Value = 1.040608852803667
$ws.Range("I17").Value = 1.034919029931369
$ws.Range("J17").Value = 1.016509374661475
$ws.Range("K17").Value = 1.042594786794224
$ws.Range("L17").Value = 1.015962329440128
$ws.Range("M17").Value = 1.043958709475527
$ws.Range("N17").Value = 1.009468525208952
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.010501771091669
$ws.Range("D18").Value = 1.039414930256067
$ws.Range("E18").Value = 1.012815214127875
$ws.Range("F18").Value = 1.040864524679987
$ws.Range("I18").Value = 1.034964087447611
$ws.Range("J18").Value = 1.016778299643135
$ws.Range("K18").Value = 1.042727163286097
$ws.Range("L18").Value = 1.016220893158166
$ws.Range("M18").Value = 1.044171816165937
$ws.Range("N18").Value = 1.009560340867985
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.010622032776173
$ws.Range("D19").Value = 1.039474451715407
$ws.Range("E19").Value = 1.012918391846671
$ws.Range("F19").Value = 1.040951641019611
$ws.Range("I19").Value = 1.034979389499085
$ws.Range("J19").Value = 1.016869952841779
$ws.Range("K19").Value = 1.042772230025626
$ws.Range("L19").Value = 1.016309019528355
$ws.Range("M19").Value = 1.04424440748918
$ws.Range("N19").Value = 1.009591627005474
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.010084040766291
$ws.Range("D20").Value = 1.039208032441194
$ws.Range("E20").Value = 1.012456850118124
$ws.Range("F20").Value = 1.040561794831594
$ws.Range("I20").Value = 1.034910712797152
$ws.Range("J20").Value = 1.016459887343085
$ws.Range("K20").Value = 1.042570403796877
$ws.Range("L20").Value = 1.015914750894257
$ws.Range("M20").Value = 1.043919475647241
$ws.Range("N20").Value = 1.009451626599729
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.00833223825645
$ws.Range("D21").Value = 1.038337918111481
$ws.Range("E21").Value = 1.010954427599715
$ws.Range("F21").Value = 1.039290114744875
$ws.Range("I21").Value = 1.034683182889366
$ws.Range("J21").Value = 1.015123694042889
$ws.Range("K21").Value = 1.041909377107472
$ws.Range("L21").Value = 1.014630339020437
$ws.Range("M21").Value = 1.042858051249092
$ws.Range("N21").Value = 1.008995031398026
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.0072286572832
$ws.Range("D22").Value = 1.037787802655467
$ws.Range("E22").Value = 1.010008294312591
$ws.Range("F22").Value = 1.03848728223649
$ws.Range("I22").Value = 1.034536862829734
$ws.Range("J22").Value = 1.014281214565481
$ws.Range("K22").Value = 1.041490025263053
$ws.Range("L22").Value = 1.013820744166228
$ws.Range("M22").Value = 1.042186813850515
$ws.Range("N22").Value = 1.008706834860351
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.007813894056029
$ws.Range("D23").Value = 1.038079717765141
$ws.Range("E23").Value = 1.010510002302008
$ws.Range("F23").Value = 1.038913190667339
$ws.Range("I23").Value = 1.034614739240755
$ws.Range("J23").Value = 1.014728055253068
$ws.Range("K23").Value = 1.041712686285216
$ws.Range("L23").Value = 1.01425012088774
$ws.Range("M23").Value = 1.042543018037786
$ws.Range("N23").Value = 1.008859720021339
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.010113374250143
$ws.Range("D24").Value = 1.039222568563047
$ws.Range("E24").Value = 1.012482013569552
$ws.Range("F24").Value = 1.040583059420102
$ws.Range("I24").Value = 1.034914472070906
$ws.Range("J24").Value = 1.016482249324647
$ws.Range("K24").Value = 1.042581422703999
$ws.Range("L24").Value = 1.0159362502724
$ws.Range("M24").Value = 1.043937205049115
$ws.Range("N24").Value = 1.009459262731106
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.012772767948574
$ws.Range("D25").Value = 1.040535577693279
$ws.Range("E25").Value = 1.014764146761313
$ws.Range("F25").Value = 1.042506691823362
$ws.Range("I25").Value = 1.035248094759934
$ws.Range("J25").Value = 1.018507860394927
$ws.Range("K25").Value = 1.043573298397437
$ws.Range("L25").Value = 1.017884287207982
$ws.Range("M25").Value = 1.04553829198262
$ws.Range("N25").Value = 1.010150218813517
